# Apply weekly data refresh to "Hortaliza, Macroferia Regional de Talca - Pepino dulce"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44245
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 667

# Row 3
$ws.Range("D3").Value = 44245
$ws.Range("I3").Value = 'Segunda'
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 556

# Row 4
$ws.Range("D4").Value = 44391
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 15000
$ws.Range("P4").Value = 833

# Row 5
$ws.Range("D5").Value = 44630
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("P5").Value = 833

# Row 6
$ws.Range("D6").Value = 44235
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 13000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 13000
$ws.Range("P6").Value = 722

# Row 7
$ws.Range("D7").Value = 44235
$ws.Range("K7").Value = 11000
$ws.Range("L7").Value = 11000
$ws.Range("M7").Value = 11000
$ws.Range("P7").Value = 611

# Row 8
$ws.Range("D8").Value = 44235
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 9000
$ws.Range("P8").Value = 500

# Row 9
$ws.Range("D9").Value = 44614
$ws.Range("J9").Value = 300
$ws.Range("N9").Value = '$/caja 18 kilos granel'

# Row 10
$ws.Range("D10").Value = 44229
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 15000
$ws.Range("P10").Value = 833

# Row 11
$ws.Range("D11").Value = 44596
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 14000
$ws.Range("P11").Value = 778

# Row 12
$ws.Range("D12").Value = 44383
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 16000
$ws.Range("P12").Value = 889

# Row 13
$ws.Range("D13").Value = 44383
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("N13").Value = '$/bandeja 18 kilos'
$ws.Range("P13").Value = 667

# Row 14
$ws.Range("D14").Value = 44396
$ws.Range("J14").Value = 250

# Row 15
$ws.Range("D15").Value = 44396
$ws.Range("I15").Value = 'Segunda'
$ws.Range("K15").Value = 12000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 12000
$ws.Range("P15").Value = 667

# Row 16
$ws.Range("D16").Value = 44249
$ws.Range("J16").Value = 400

# Row 17
$ws.Range("D17").Value = 44249

# Row 18
$ws.Range("D18").Value = 44628
$ws.Range("J18").Value = 300

# Row 19
$ws.Range("D19").Value = 44585
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 200

# Row 20
$ws.Range("D20").Value = 44238
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 12000
$ws.Range("P20").Value = 667

# Row 21
$ws.Range("D21").Value = 44238
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 10000
$ws.Range("P21").Value = 556

# Row 22
$ws.Range("D22").Value = 44238
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 8000
$ws.Range("P22").Value = 444

# Row 23
$ws.Range("D23").Value = 44635

# Row 24
$ws.Range("D24").Value = 44627
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("P24").Value = 833

# Row 25
$ws.Range("D25").Value = 44243
$ws.Range("I25").Value = 'Especial'
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 12000
$ws.Range("M25").Value = 12000
$ws.Range("P25").Value = 667

# Row 26
$ws.Range("D26").Value = 44243
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 10000
$ws.Range("P26").Value = 556

# Row 27
$ws.Range("D27").Value = 44243
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 8000
$ws.Range("P27").Value = 444
